$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2444.4443
$ws.Range("I51").Value = 10001
$ws.Range("J51").Value = 1499.875
$ws.Range("K51").Value = 10001
$ws.Range("L51").Value = 1499.875
$ws.Range("M51").Value = -9517
$ws.Range("N51").Value = -2467.875
$ws.Range("H64").Value = 3511.5
$ws.Range("I64").Value = 3349.5186
$ws.Range("K64").Value = 3349.5186
$ws.Range("M64").Value = -3101.5186
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3511.5
$ws.Range("I67").Value = 3349.5186
$ws.Range("K67").Value = 3349.5186
$ws.Range("M67").Value = -2491.5186
$ws.Range("N67").ClearContents()
$ws.Range("H76").Value = 3728
$ws.Range("I76").Value = 3606.8276
$ws.Range("K76").Value = 3606.8276
$ws.Range("M76").Value = -3291.8276
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3728
$ws.Range("I79").Value = 3606.8276
$ws.Range("K79").Value = 3606.8276
$ws.Range("M79").Value = -2514.8276
$ws.Range("N79").ClearContents()
$ws.Range("H113").Value = 2837.625
$ws.Range("J113").Value = 2934.8
$ws.Range("L113").Value = 2934.8
$ws.Range("N113").Value = -9442.799999999999
$ws.Range("H138").Value = 3438.632
$ws.Range("I138").Value = 1897.7391
$ws.Range("J138").Value = 3992.3906
$ws.Range("K138").Value = 5693.2173
$ws.Range("L138").Value = 11977.1718
$ws.Range("M138").Value = -553.2173000000003
$ws.Range("N138").Value = -22257.1718

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6481.5386
$ws.Range("I61").Value = 5448.1875
$ws.Range("J61").Value = 8134.9
$ws.Range("K61").Value = 5448.1875
$ws.Range("L61").Value = 8134.9
$ws.Range("M61").Value = -5236.1875
$ws.Range("N61").Value = -8558.9
$ws.Range("H74").Value = 6201.6
$ws.Range("I74").Value = 3349.7144
$ws.Range("J74").Value = 12856
$ws.Range("K74").Value = 3349.7144
$ws.Range("L74").Value = 12856
$ws.Range("M74").Value = -2475.7144
$ws.Range("N74").Value = -14604
$ws.Range("H77").Value = 6201.6
$ws.Range("I77").Value = 3349.7144
$ws.Range("J77").Value = 12856
$ws.Range("K77").Value = 16748.572
$ws.Range("L77").Value = 64280
$ws.Range("M77").Value = -12380.572
$ws.Range("N77").Value = -73016
$ws.Range("H122").Value = 5437663.5
$ws.Range("I122").Value = 3387.3845
$ws.Range("J122").Value = 12502223
$ws.Range("K122").Value = 10162.1535
$ws.Range("L122").Value = 37506669
$ws.Range("M122").Value = -7712.1535
$ws.Range("N122").Value = -37511569
$ws.Range("H136").Value = 6481.5386
$ws.Range("I136").Value = 5448.1875
$ws.Range("J136").Value = 8134.9
$ws.Range("K136").Value = 16344.5625
$ws.Range("L136").Value = 24404.7
$ws.Range("M136").Value = -13794.5625
$ws.Range("N136").Value = -29504.7
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 8793.111000000001
$ws.Range("I75").Value = 5892.25
$ws.Range("K75").Value = 5892.25
$ws.Range("M75").Value = -4956.25
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 8793.111000000001
$ws.Range("I78").Value = 5892.25
$ws.Range("K78").Value = 17676.75
$ws.Range("M78").Value = -12996.75
$ws.Range("N78").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 56324.75
$ws.Range("J111").Value = 56324.75
$ws.Range("L111").Value = 56324.75
$ws.Range("N111").Value = -64504.75
$ws.Range("H134").Value = 4566.077
$ws.Range("I134").Value = 4566.077
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13698.231
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11163.231
$ws.Range("N134").ClearContents()
$ws.Range("H141").Value = 43022.715
$ws.Range("J141").Value = 43022.715
$ws.Range("L141").Value = 43022.715
$ws.Range("N141").Value = -53382.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 82.86667
$ws.Range("I7").Value = 81.666664
$ws.Range("J7").Value = 84.666664
$ws.Range("K7").Value = 81.666664
$ws.Range("L7").Value = 84.666664
$ws.Range("M7").Value = 31.333336
$ws.Range("N7").Value = -310.666664
$ws.Range("H41").Value = 25860
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9572
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value = 6997077
$ws.Range("I58").Value = 11367148
$ws.Range("J58").Value = 4964
$ws.Range("K58").Value = 11367148
$ws.Range("L58").Value = 4964
$ws.Range("M58").Value = -11366945
$ws.Range("N58").Value = -5370
$ws.Range("H59").Value = 41563.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 41563.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 41563.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -43853.5
$ws.Range("H68").Value = 25147.5
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41793
$ws.Range("H71").Value = 25147.5
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128373
$ws.Range("H117").Value = 60700
$ws.Range("J117").Value = 63625
$ws.Range("L117").Value = 63625
$ws.Range("N117").Value = -72803
$ws.Range("H122").Value = 12453.529
$ws.Range("I122").Value = 4526.909
$ws.Range("J122").Value = 26985.666
$ws.Range("K122").Value = 13580.727
$ws.Range("L122").Value = 80956.99800000001
$ws.Range("M122").Value = -11130.727
$ws.Range("N122").Value = -85856.99800000001
$ws.Range("H132").Value = 2980.2727
$ws.Range("I132").Value = 2492.1765
$ws.Range("J132").Value = 4639.8
$ws.Range("K132").Value = 7476.529500000001
$ws.Range("L132").Value = 13919.4
$ws.Range("M132").Value = -4946.529500000001
$ws.Range("N132").Value = -18979.4
$ws.Range("H134").Value = 4645.029
$ws.Range("I134").Value = 4011.3
$ws.Range("J134").Value = 4898.52
$ws.Range("K134").Value = 12033.9
$ws.Range("L134").Value = 14695.56
$ws.Range("M134").Value = -9498.900000000001
$ws.Range("N134").Value = -19765.56
$ws.Range("H136").Value = 6997077
$ws.Range("I136").Value = 11367148
$ws.Range("J136").Value = 4964
$ws.Range("K136").Value = 34101444
$ws.Range("L136").Value = 14892
$ws.Range("M136").Value = -34098894
$ws.Range("N136").Value = -19992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 14289187
$ws.Range("I5").Value = 436.75
$ws.Range("J5").Value = 45464640
$ws.Range("K5").Value = 1310.25
$ws.Range("L5").Value = 136393920
$ws.Range("M5").Value = -1198.25
$ws.Range("N5").Value = -136394144
$ws.Range("H93").Value = 4562.4062
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 4562.4062
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 13687.2186
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -17431.2186
$ws.Range("H132").Value = 2465.3845
$ws.Range("I132").Value = 5400
$ws.Range("J132").Value = 1931.8182
$ws.Range("K132").Value = 48600
$ws.Range("L132").Value = 17386.3638
$ws.Range("M132").Value = -46070
$ws.Range("N132").Value = -22446.3638
$ws.Range("H135").Value = 14289187
$ws.Range("I135").Value = 436.75
$ws.Range("J135").Value = 45464640
$ws.Range("K135").Value = 3930.75
$ws.Range("L135").Value = 409181760
$ws.Range("M135").Value = -1395.75
$ws.Range("N135").Value = -409186830
$ws.Range("H139").Value = 1807828.5
$ws.Range("I139").Value = 2516811.8
$ws.Range("J139").Value = 3143.9092
$ws.Range("K139").Value = 7550435.399999999
$ws.Range("L139").Value = 9431.7276
$ws.Range("M139").Value = -7545295.399999999
$ws.Range("N139").Value = -19711.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6939.385
$ws.Range("I80").Value = 12360
$ws.Range("J80").Value = 3551.5
$ws.Range("K80").Value = 12360
$ws.Range("L80").Value = 3551.5
$ws.Range("M80").Value = -11362
$ws.Range("N80").Value = -5547.5
$ws.Range("H83").Value = 6939.385
$ws.Range("I83").Value = 12360
$ws.Range("J83").Value = 3551.5
$ws.Range("K83").Value = 61800
$ws.Range("L83").Value = 17757.5
$ws.Range("M83").Value = -56808
$ws.Range("N83").Value = -27741.5
$ws.Range("H132").Value = 7028.4287
$ws.Range("I132").Value = 7840.2
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 23520.6
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -20990.6
$ws.Range("N132").Value = -20057
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4915.974
$ws.Range("I122").Value = 4281.846
$ws.Range("J122").Value = 6184.231
$ws.Range("K122").Value = 12845.538
$ws.Range("L122").Value = 18552.693
$ws.Range("M122").Value = -10395.538
$ws.Range("N122").Value = -23452.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1840.0667
$ws.Range("I81").Value = 1460.2
$ws.Range("J81").Value = 2030
$ws.Range("K81").Value = 2920.4
$ws.Range("L81").Value = 4060
$ws.Range("M81").Value = -1859.4
$ws.Range("N81").Value = -6182
$ws.Range("H84").Value = 1840.0667
$ws.Range("I84").Value = 1460.2
$ws.Range("J84").Value = 2030
$ws.Range("K84").Value = 14602
$ws.Range("L84").Value = 20300
$ws.Range("M84").Value = -9298
$ws.Range("N84").Value = -30908
$ws.Range("H136").Value = 6315.6343
$ws.Range("I136").Value = 2944
$ws.Range("J136").Value = 10623.833
$ws.Range("K136").Value = 8832
$ws.Range("L136").Value = 31871.499
$ws.Range("M136").Value = -6282
$ws.Range("N136").Value = -36971.499
